$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column for the rows being updated so that
# numeric-looking strings (e.g. "1.09", "0.999") are stored as text, matching
# the original inlineStr cell type, instead of being coerced to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.862.01"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "3.627.34"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "586.45"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").Value = "195.27"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").Value = "3.621.16"
$ws.Range("E7").Value = "  -1.20%  "

$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "0.685"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").Value = "55.56"
$ws.Range("E12").Value = "  -3.19%  "

$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  +3.47%  "

$ws.Range("D14").Value = "10.06"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "4.187.66"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "3.618.75"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "12.49"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "18.58"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("D20").Value = "67.765.29"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("D21").Value = "1.09"
$ws.Range("E21").Value = "  -2.43%  "

$ws.Range("D22").Value = "405.72"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").Value = "13.41"
$ws.Range("E23").Value = "  +19.89%  "

$ws.Range("D24").Value = "4.28"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("D25").Value = "86.33"
$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("D26").Value = "2.96"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "12.71"
$ws.Range("E27").Value = "  +0.46%  "

$ws.Range("D28").Value = "3.88"
$ws.Range("E28").Value = "  +5.04%  "

$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  +12.37%  "

$ws.Range("D31").Value = "9.27"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("D32").Value = "31.58"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("D33").Value = "670.78"
$ws.Range("E33").Value = "  +9.57%  "

$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "64.52"
$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "43.27"
$ws.Range("E37").Value = "  -3.98%  "

$ws.Range("D38").Value = "0.424"
$ws.Range("E38").Value = "  +7.09%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0802"
$ws.Range("E39").Value = "  +3.79%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +16.72%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.178.70"
$ws.Range("E42").Value = "  +14.49%  "

$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.11"
$ws.Range("E43").Value = "  +7.39%  "

$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").Value = "0.0422"
$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "8.89"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.132"
$ws.Range("E48").Value = "  -2.61%  "

$ws.Range("D49").Value = "3.12"
$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("D50").Value = "143.22"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").Value = "2.75"
$ws.Range("E51").Value = "  +1.68%  "

# Restore the default (Normal) style on the price column so no stray number
# format is left behind on cells now holding text values.
$priceRange.Style = "Normal"
